$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 35.73885133333334
$ws.Cells.Item(2, 8).Value = 107.216554
$ws.Cells.Item(2, 9).Value = 0.01949729408921566
$ws.Cells.Item(2, 10).Value = 0.01949729408921566
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.971099
$ws.Cells.Item(2, 14).Value = 2.913297
$ws.Cells.Item(2, 15).Value = 0.007882574716876797
$ws.Cells.Item(2, 16).Value = 0.007882574716876797
$ws.Cells.Item(2, 17).Value = 34.70596279094867
$ws.Cells.Item(2, 18).Value = 312.353665118538
$ws.Cells.Item(2, 19).Value = 0.0001536888774351628
$ws.Cells.Item(2, 20).Value = 0.0001536888774351628

# Row 3
$ws.Cells.Item(3, 7).Value = 35.73885133333334
$ws.Cells.Item(3, 8).Value = 107.216554
$ws.Cells.Item(3, 9).Value = 0.01949729408921566
$ws.Cells.Item(3, 10).Value = 0.01949729408921566
$ws.Cells.Item(3, 15).Value = 0.6966643430097871
$ws.Cells.Item(3, 16).Value = 0.696664343009787
$ws.Cells.Item(3, 17).Value = 3067.323512267099
$ws.Cells.Item(3, 18).Value = 27605.91161040389
$ws.Cells.Item(3, 19).Value = 0.01358306957713203
$ws.Cells.Item(3, 20).Value = 0.01358306957713203

# Row 4
$ws.Cells.Item(4, 7).Value = 35.73885133333334
$ws.Cells.Item(4, 8).Value = 107.216554
$ws.Cells.Item(4, 9).Value = 0.01949729408921566
$ws.Cells.Item(4, 10).Value = 0.01949729408921566
$ws.Cells.Item(4, 13).Value = 36.24916566666667
$ws.Cells.Item(4, 14).Value = 108.747497
$ws.Cells.Item(4, 15).Value = 0.294240604502677
$ws.Cells.Item(4, 16).Value = 0.294240604502677
$ws.Cells.Item(4, 17).Value = 1295.503542718371
$ws.Cells.Item(4, 18).Value = 11659.53188446534
$ws.Cells.Item(4, 19).Value = 0.005736895598977288
$ws.Cells.Item(4, 20).Value = 0.005736895598977286

# Row 5
$ws.Cells.Item(5, 7).Value = 35.73885133333334
$ws.Cells.Item(5, 8).Value = 107.216554
$ws.Cells.Item(5, 9).Value = 0.01949729408921566
$ws.Cells.Item(5, 10).Value = 0.01949729408921566
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.149372
$ws.Cells.Item(5, 14).Value = 0.448116
$ws.Cells.Item(5, 15).Value = 0.001212477770659141
$ws.Cells.Item(5, 16).Value = 0.001212477770659141
$ws.Cells.Item(5, 17).Value = 5.338383701362666
$ws.Cells.Item(5, 18).Value = 48.04545331226399
$ws.Cells.Item(5, 19).Value = 0.00002364003567117784
$ws.Cells.Item(5, 20).Value = 0.00002364003567117784

# Row 6
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.9215900675332435
$ws.Cells.Item(6, 10).Value = 0.9215900675332435
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.971099
$ws.Cells.Item(6, 14).Value = 2.913297
$ws.Cells.Item(6, 15).Value = 0.007882574716876797
$ws.Cells.Item(6, 16).Value = 0.007882574716876797
$ws.Cells.Item(6, 17).Value = 1640.467156414693
$ws.Cells.Item(6, 18).Value = 14764.20440773224
$ws.Cells.Item(6, 19).Value = 0.007264502565662325
$ws.Cells.Item(6, 20).Value = 0.007264502565662325

# Row 7
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.9215900675332435
$ws.Cells.Item(7, 10).Value = 0.9215900675332435
$ws.Cells.Item(7, 15).Value = 0.6966643430097871
$ws.Cells.Item(7, 16).Value = 0.696664343009787
$ws.Cells.Item(7, 19).Value = 0.6420389389223924
$ws.Cells.Item(7, 20).Value = 0.6420389389223923

# Row 8
$ws.Cells.Item(8, 8).Value = 5067.86792
$ws.Cells.Item(8, 9).Value = 0.9215900675332435
$ws.Cells.Item(8, 10).Value = 0.9215900675332435
$ws.Cells.Item(8, 13).Value = 36.24916566666667
$ws.Cells.Item(8, 14).Value = 108.747497
$ws.Cells.Item(8, 15).Value = 0.294240604502677
$ws.Cells.Item(8, 16).Value = 0.294240604502677
$ws.Cells.Item(8, 17).Value = 61235.32793628847
$ws.Cells.Item(8, 18).Value = 551117.9514265963
$ws.Cells.Item(8, 19).Value = 0.2711692185746445
$ws.Cells.Item(8, 20).Value = 0.2711692185746444

# Row 9
$ws.Cells.Item(9, 8).Value = 5067.86792
$ws.Cells.Item(9, 9).Value = 0.9215900675332435
$ws.Cells.Item(9, 10).Value = 0.9215900675332435
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.149372
$ws.Cells.Item(9, 14).Value = 0.448116
$ws.Cells.Item(9, 15).Value = 0.001212477770659141
$ws.Cells.Item(9, 16).Value = 0.001212477770659141
$ws.Cells.Item(9, 17).Value = 252.3325223154133
$ws.Cells.Item(9, 18).Value = 2270.992700838719
$ws.Cells.Item(9, 19).Value = 0.001117407470544314
$ws.Cells.Item(9, 20).Value = 0.001117407470544314

# Row 10
$ws.Cells.Item(10, 7).Value = 93.641553
$ws.Cells.Item(10, 8).Value = 280.924659
$ws.Cells.Item(10, 9).Value = 0.05108605424341119
$ws.Cells.Item(10, 10).Value = 0.05108605424341119
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.971099
$ws.Cells.Item(10, 14).Value = 2.913297
$ws.Cells.Item(10, 15).Value = 0.007882574716876797
$ws.Cells.Item(10, 16).Value = 0.007882574716876797
$ws.Cells.Item(10, 17).Value = 90.93521847674701
$ws.Cells.Item(10, 18).Value = 818.416966290723
$ws.Cells.Item(10, 19).Value = 0.0004026896395641096
$ws.Cells.Item(10, 20).Value = 0.0004026896395641096

# Row 11
$ws.Cells.Item(11, 7).Value = 93.641553
$ws.Cells.Item(11, 8).Value = 280.924659
$ws.Cells.Item(11, 9).Value = 0.05108605424341119
$ws.Cells.Item(11, 10).Value = 0.05108605424341119
$ws.Cells.Item(11, 15).Value = 0.6966643430097871
$ws.Cells.Item(11, 16).Value = 0.696664343009787
$ws.Cells.Item(11, 17).Value = 8036.882175175272
$ws.Cells.Item(11, 18).Value = 72331.93957657745
$ws.Cells.Item(11, 19).Value = 0.0355898324164484
$ws.Cells.Item(11, 20).Value = 0.03558983241644839

# Row 12
$ws.Cells.Item(12, 7).Value = 93.641553
$ws.Cells.Item(12, 8).Value = 280.924659
$ws.Cells.Item(12, 9).Value = 0.05108605424341119
$ws.Cells.Item(12, 10).Value = 0.05108605424341119
$ws.Cells.Item(12, 13).Value = 36.24916566666667
$ws.Cells.Item(12, 14).Value = 108.747497
$ws.Cells.Item(12, 15).Value = 0.294240604502677
$ws.Cells.Item(12, 16).Value = 0.294240604502677
$ws.Cells.Item(12, 17).Value = 3394.428167980947
$ws.Cells.Item(12, 18).Value = 30549.85351182853
$ws.Cells.Item(12, 19).Value = 0.01503159148223786
$ws.Cells.Item(12, 20).Value = 0.01503159148223785

# Row 13
$ws.Cells.Item(13, 7).Value = 93.641553
$ws.Cells.Item(13, 8).Value = 280.924659
$ws.Cells.Item(13, 9).Value = 0.05108605424341119
$ws.Cells.Item(13, 10).Value = 0.05108605424341119
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.149372
$ws.Cells.Item(13, 14).Value = 0.448116
$ws.Cells.Item(13, 15).Value = 0.001212477770659141
$ws.Cells.Item(13, 16).Value = 0.001212477770659141
$ws.Cells.Item(13, 17).Value = 13.987426054716
$ws.Cells.Item(13, 18).Value = 125.886834492444
$ws.Cells.Item(13, 19).Value = 0.00006194070516082313
$ws.Cells.Item(13, 20).Value = 0.00006194070516082313

# Row 14
$ws.Cells.Item(14, 7).Value = 14.34625366666667
$ws.Cells.Item(14, 8).Value = 43.038761
$ws.Cells.Item(14, 9).Value = 0.007826584134129748
$ws.Cells.Item(14, 10).Value = 0.007826584134129748
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.971099
$ws.Cells.Item(14, 14).Value = 2.913297
$ws.Cells.Item(14, 15).Value = 0.007882574716876797
$ws.Cells.Item(14, 16).Value = 0.007882574716876797
$ws.Cells.Item(14, 17).Value = 13.93163258944633
$ws.Cells.Item(14, 18).Value = 125.384693305017
$ws.Cells.Item(14, 19).Value = 0.00006169363421520023
$ws.Cells.Item(14, 20).Value = 0.00006169363421520023

# Row 15
$ws.Cells.Item(15, 7).Value = 14.34625366666667
$ws.Cells.Item(15, 8).Value = 43.038761
$ws.Cells.Item(15, 9).Value = 0.007826584134129748
$ws.Cells.Item(15, 10).Value = 0.007826584134129748
$ws.Cells.Item(15, 15).Value = 0.6966643430097871
$ws.Cells.Item(15, 16).Value = 0.696664343009787
$ws.Cells.Item(15, 17).Value = 1231.281911505421
$ws.Cells.Item(15, 18).Value = 11081.53720354879
$ws.Cells.Item(15, 19).Value = 0.005452502093814324
$ws.Cells.Item(15, 20).Value = 0.005452502093814323

# Row 16
$ws.Cells.Item(16, 7).Value = 14.34625366666667
$ws.Cells.Item(16, 8).Value = 43.038761
$ws.Cells.Item(16, 9).Value = 0.007826584134129748
$ws.Cells.Item(16, 10).Value = 0.007826584134129748
$ws.Cells.Item(16, 13).Value = 36.24916566666667
$ws.Cells.Item(16, 14).Value = 108.747497
$ws.Cells.Item(16, 15).Value = 0.294240604502677
$ws.Cells.Item(16, 16).Value = 0.294240604502677
$ws.Cells.Item(16, 17).Value = 520.0397258590242
$ws.Cells.Item(16, 18).Value = 4680.357532731217
$ws.Cells.Item(16, 19).Value = 0.002302898846817398
$ws.Cells.Item(16, 20).Value = 0.002302898846817398

# Row 17
$ws.Cells.Item(17, 7).Value = 14.34625366666667
$ws.Cells.Item(17, 8).Value = 43.038761
$ws.Cells.Item(17, 9).Value = 0.007826584134129748
$ws.Cells.Item(17, 10).Value = 0.007826584134129748
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.149372
$ws.Cells.Item(17, 14).Value = 0.448116
$ws.Cells.Item(17, 15).Value = 0.001212477770659141
$ws.Cells.Item(17, 16).Value = 0.001212477770659141
$ws.Cells.Item(17, 17).Value = 2.142928602697333
$ws.Cells.Item(17, 18).Value = 19.286357424276
$ws.Cells.Item(17, 19).Value = 0.000009489559282825839
$ws.Cells.Item(17, 20).Value = 0.000009489559282825839
